$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.33622468515676
$ws.Range("C2").Value = 4.24680420599453
$ws.Range("D2").Value = 9.042552701705485
$ws.Range("E2").Value = 13.63193449268419
$ws.Range("F2").Value = 33.75237425238932
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("J2").Value = 9.969556080289117
$ws.Range("K2").Value = 10.71800481330001
$ws.Range("N2").Value = 19.32907680124845
$ws.Range("O2").Value = 25.54615279963842

$ws.Range("B3").Value = 11.08334338012602
$ws.Range("C3").Value = 4.023729007910274
$ws.Range("D3").Value = 8.982006661439817
$ws.Range("E3").Value = 13.56835971976752
$ws.Range("F3").Value = 33.79094773599694
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("J3").Value = 9.975706601826554
$ws.Range("K3").Value = 10.54684682131505
$ws.Range("N3").Value = 19.39007248772496
$ws.Range("O3").Value = 25.61521826754151

$ws.Range("B4").Value = 10.92732458837222
$ws.Range("C4").Value = 3.879456251352135
$ws.Range("D4").Value = 8.946288665051929
$ws.Range("E4").Value = 13.53206170138489
$ws.Range("F4").Value = 33.82232951032461
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("J4").Value = 9.981048082998315
$ws.Range("K4").Value = 10.44218786508059
$ws.Range("N4").Value = 19.42926453858485
$ws.Range("O4").Value = 25.66271270910531

$ws.Range("B5").Value = 10.8636521272036
$ws.Range("C5").Value = 3.818859765321839
$ws.Range("D5").Value = 8.932112424210437
$ws.Range("E5").Value = 13.51797013938423
$ws.Range("F5").Value = 33.83705123336836
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("J5").Value = 9.983618710367329
$ws.Range("K5").Value = 10.39970067960932
$ws.Range("N5").Value = 19.44567458088816
$ws.Range("O5").Value = 25.68334431743715

$ws.Range("B6").Value = 10.85307654011704
$ws.Range("C6").Value = 3.808689879186618
$ws.Range("D6").Value = 8.929781726051905
$ws.Range("E6").Value = 13.51567286425034
$ws.Range("F6").Value = 33.83961245800076
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("J6").Value = 9.984069362787915
$ws.Range("K6").Value = 10.39265707545961
$ws.Range("N6").Value = 19.4484260068301
$ws.Range("O6").Value = 25.6868472562441

$ws.Range("B7").Value = 10.92646613068221
$ws.Range("C7").Value = 3.878646283178275
$ws.Range("D7").Value = 8.946095927976961
$ws.Range("E7").Value = 13.53186880764541
$ws.Range("F7").Value = 33.82252022785546
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("J7").Value = 9.981081155932685
$ws.Range("K7").Value = 10.44161414063487
$ws.Range("N7").Value = 19.42948407107261
$ws.Range("O7").Value = 25.66298578568843

$ws.Range("B8").Value = 11.24924244445615
$ws.Range("C8").Value = 4.171421182082017
$ws.Range("D8").Value = 9.021381674150639
$ws.Range("E8").Value = 13.60945323379126
$ws.Range("F8").Value = 33.76407565702562
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("J8").Value = 9.971352220041325
$ws.Range("K8").Value = 10.65893162965605
$ws.Range("N8").Value = 19.34974767665815
$ws.Range("O8").Value = 25.56890929775303

$ws.Range("B9").Value = 11.87227593528205
$ws.Range("C9").Value = 4.686488946344082
$ws.Range("D9").Value = 9.180000155021526
$ws.Range("E9").Value = 13.7827492692501
$ws.Range("F9").Value = 33.71061283972622
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("J9").Value = 9.964671247501002
$ws.Range("K9").Value = 11.08620502099724
$ws.Range("N9").Value = 19.20713362415497
$ws.Range("O9").Value = 25.42488916265347

$ws.Range("B10").Value = 12.31890842837108
$ws.Range("C10").Value = 5.027773412406093
$ws.Range("D10").Value = 9.302420473098733
$ws.Range("E10").Value = 13.9221515361449
$ws.Range("F10").Value = 33.70866542213513
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("J10").Value = 9.967288165549398
$ws.Range("K10").Value = 11.39776038798466
$ws.Range("N10").Value = 19.11064990517346
$ws.Range("O10").Value = 25.34386567196773

$ws.Range("B11").Value = 12.51870349243874
$ws.Range("C11").Value = 5.17481213832419
$ws.Range("D11").Value = 9.359208833258014
$ws.Range("E11").Value = 13.9880030519287
$ws.Range("F11").Value = 33.71588117223126
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("J11").Value = 9.970103536249081
$ws.Range("K11").Value = 11.53836425012388
$ws.Range("N11").Value = 19.06853957198299
$ws.Range("O11").Value = 25.31241129402905

$ws.Range("B12").Value = 12.59380008679454
$ws.Range("C12").Value = 5.229301613752214
$ws.Range("D12").Value = 9.380855407790515
$ws.Range("E12").Value = 14.01327268560497
$ws.Range("F12").Value = 33.71977668883905
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("J12").Value = 9.971402296710091
$ws.Range("K12").Value = 11.59139645286061
$ws.Range("N12").Value = 19.05284816355173
$ws.Range("O12").Value = 25.3012789670536

$ws.Range("B13").Value = 12.57765280612815
$ws.Range("C13").Value = 5.217619425923602
$ws.Range("D13").Value = 9.376187368526114
$ws.Range("E13").Value = 14.00781588503489
$ws.Range("F13").Value = 33.71888603194335
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("J13").Value = 9.971112256236104
$ws.Range("K13").Value = 11.57998519228792
$ws.Range("N13").Value = 19.05621627288575
$ws.Range("O13").Value = 25.30364185155911

$ws.Range("B14").Value = 12.52489339975292
$ws.Range("C14").Value = 5.179318935079358
$ws.Range("D14").Value = 9.36098695495801
$ws.Range("E14").Value = 13.9900754253609
$ws.Range("F14").Value = 33.71617836328579
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("J14").Value = 9.970205729879554
$ws.Range("K14").Value = 11.54273175622753
$ws.Range("N14").Value = 19.06724352933519
$ws.Range("O14").Value = 25.31147981395905

$ws.Range("B15").Value = 12.49250147512428
$ws.Range("C15").Value = 5.155703439499151
$ws.Range("D15").Value = 9.351694286106676
$ws.Range("E15").Value = 13.97925173678215
$ws.Range("F15").Value = 33.71467123159256
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("J15").Value = 9.969680720564128
$ws.Range("K15").Value = 11.51988397577942
$ws.Range("N15").Value = 19.07403119768933
$ws.Range("O15").Value = 25.3163822584167

$ws.Range("B16").Value = 12.30577676586264
$ws.Range("C16").Value = 5.01799791354272
$ws.Range("D16").Value = 9.29872997146445
$ws.Range("E16").Value = 13.9178955521163
$ws.Range("F16").Value = 33.70835673581833
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("J16").Value = 9.96713678759202
$ws.Range("K16").Value = 11.38854467088103
$ws.Range("N16").Value = 19.11343764292373
$ws.Range("O16").Value = 25.34603017615582

$ws.Range("B17").Value = 12.19030611210327
$ws.Range("C17").Value = 4.931408578331683
$ws.Range("D17").Value = 9.266508217620974
$ws.Range("E17").Value = 13.88086750237046
$ws.Range("F17").Value = 33.70655691845108
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("J17").Value = 9.965991661571422
$ws.Range("K17").Value = 11.30764825626559
$ws.Range("N17").Value = 19.13806742604474
$ws.Range("O17").Value = 25.36560349759319

$ws.Range("B18").Value = 12.12357644907461
$ws.Range("C18").Value = 4.880832015915662
$ws.Range("D18").Value = 9.248079505330779
$ws.Range("E18").Value = 13.85980037940445
$ws.Range("F18").Value = 33.70628444649622
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("J18").Value = 9.965486029662593
$ws.Range("K18").Value = 11.26101581099818
$ws.Range("N18").Value = 19.15240149405546
$ws.Range("O18").Value = 25.3773700759965

$ws.Range("B19").Value = 12.10093138632477
$ws.Range("C19").Value = 4.863575297050705
$ws.Range("D19").Value = 9.241858271195598
$ws.Range("E19").Value = 13.85270750794805
$ws.Range("F19").Value = 33.70632323756514
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("J19").Value = 9.965341140808837
$ws.Range("K19").Value = 11.24521073180788
$ws.Range("N19").Value = 19.15728359670773
$ws.Range("O19").Value = 25.38144132066009

$ws.Range("B20").Value = 12.20263122074559
$ws.Range("C20").Value = 4.940706208365874
$ws.Range("D20").Value = 9.269927589393687
$ws.Range("E20").Value = 13.88478547613419
$ws.Range("F20").Value = 33.70666958309027
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("J20").Value = 9.96609773278851
$ws.Range("K20").Value = 11.3162708415275
$ws.Range("N20").Value = 19.1354281991746
$ws.Range("O20").Value = 25.36346724090081

$ws.Range("B21").Value = 12.54040592846647
$ws.Range("C21").Value = 5.190601104266064
$ws.Range("D21").Value = 9.365447957532812
$ws.Range("E21").Value = 13.99527732908542
$ws.Range("F21").Value = 33.71694212560473
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("J21").Value = 9.970465693599245
$ws.Range("K21").Value = 11.55368011765276
$ws.Range("N21").Value = 19.06399764936356
$ws.Range("O21").Value = 25.3091564661215

$ws.Range("B22").Value = 12.75785740905395
$ws.Range("C22").Value = 5.346977516694239
$ws.Range("D22").Value = 9.428696535029388
$ws.Range("E22").Value = 14.06942325355137
$ws.Range("F22").Value = 33.73043356641631
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("J22").Value = 9.974675854710595
$ws.Range("K22").Value = 11.7075868806754
$ws.Range("N22").Value = 19.01879863684462
$ws.Range("O22").Value = 25.27820113454456

$ws.Range("B23").Value = 12.64212528444784
$ws.Range("C23").Value = 5.264154510521908
$ws.Range("D23").Value = 9.394869814903201
$ws.Range("E23").Value = 14.02967915956223
$ws.Range("F23").Value = 33.72261363180682
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("J23").Value = 9.972305158607853
$ws.Range("K23").Value = 11.62557456533165
$ws.Range("N23").Value = 19.04278671384787
$ws.Range("O23").Value = 25.29430665295631

$ws.Range("B24").Value = 12.19706010495574
$ws.Range("C24").Value = 4.9365052259562
$ws.Range("D24").Value = 9.268381391174898
$ws.Range("E24").Value = 13.88301347114561
$ws.Range("F24").Value = 33.7066162728442
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("J24").Value = 9.966049302221396
$ws.Range("K24").Value = 11.31237295320034
$ws.Range("N24").Value = 19.13662085129461
$ws.Range("O24").Value = 25.3644314427843

$ws.Range("B25").Value = 11.70533160283837
$ws.Range("C25").Value = 4.553606133814391
$ws.Range("D25").Value = 9.135995941296491
$ws.Range("E25").Value = 13.73368429218419
$ws.Range("F25").Value = 33.71852090139635
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("J25").Value = 9.965154271169785
$ws.Range("K25").Value = 10.97081581670394
$ws.Range("N25").Value = 19.24425161500143
$ws.Range("O25").Value = 25.45950547494647
